# Auto commit at 2025-08-29  7:46:45.98
#
# Updates the monthly "Metrics" figures (month-of-August refresh) and lets
# the dependent "today" sheet (which pulls every number from Metrics! via
# formulas, plus a TODAY()-1 "as of" date) recompute on its own through
# normal Excel automatic recalculation.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsToday   = $wb.Worksheets.Item("today")

# --- Metrics!B2:B13 -> refreshed totals -------------------------------
$wsMetrics.Range("B2").Value  = 473093.14
$wsMetrics.Range("B3").Value  = 405842.87
$wsMetrics.Range("B4").Value  = 149681.38999999998
$wsMetrics.Range("B5").Value  = 18543
$wsMetrics.Range("B6").Value  = 3868721.7099999995
$wsMetrics.Range("B7").Value  = 3284557.53
$wsMetrics.Range("B8").Value  = 1113323.95
$wsMetrics.Range("B9").Value  = 149231
$wsMetrics.Range("B10").Value = 32334045.510999829
$wsMetrics.Range("B11").Value = 19314427.600000001
$wsMetrics.Range("B12").Value = 11395032.840000002
$wsMetrics.Range("B13").Value = 1246858

# The "today" sheet's B11:B22 (=Metrics!B2..B13), E11:E22 (=B col) and
# F11:F22 (=E col + another Metrics cell) are all formulas, so they pick
# up the refreshed Metrics values automatically on recalculation, as does
# A1's "=TODAY()-1" as-of-date cell.

# --- Restore the saved cursor position on each sheet -------------------
# Select on "today" first, then re-activate "Metrics" last so Metrics
# ends up the active (tabSelected) sheet, matching the saved view state.
$wsToday.Activate()
$wsToday.Range("I9").Select() | Out-Null

$wsMetrics.Activate()
$wsMetrics.Range("E19").Select() | Out-Null
